$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the "Latest HO Xliff Generate Date" column on the Overview sheet
# for the four files that just finished handoff generation.
$overview.Range("G4:G7").Value = "2016-08-24 08:32:35"

# zh-cn sheet: these four rows are now handed-off ("ht" priority) with a
# fresh Latest Handoff Datetime.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-24 08:32:30"

# de-de sheet: same priority bump, plus the shared "Latest Handoff Datetime"
# entry (which coincides with the Overview generate-date string).
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-24 08:32:35"
